# Apply the field updates described by the diff.
$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $true, $false, $false, $false, `
                             $true, 1, $false, $replace, 2)
}

Replace-Text "2A0508A" "2S0101A"
Replace-Text "OFICIAL ADMINISTRATIVO 5A" "AUXILIAR 1A"
Replace-Text "CONTABILIDAD GUBERNAMENTAL" "CONTROL DE REC. HUMANOS Y SUELDOS APLICADOS"
Replace-Text "1140031490300000120" "1140020000000000220"
